# The sheet currently has its data table in columns B:F (rows 1-3), with
# column A holding a couple of stray style-only values (row 2/3). The table
# needs to move one column to the left (A:E) and the "MODEL_CONDITION"
# header needs to lose its underscore ("MODELCONDITION").
#
# Copying B1:F3 onto A1:E3 shifts every value (and its cell formatting,
# e.g. the bold/border header style) left by one column in a single step,
# which is exactly what's needed: the former A2/A3 values are overwritten
# and the former F column values (2 and 9) land in the new E column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1:F3").Copy($ws.Range("A1:E3"))

# Column F is now a leftover duplicate of column E; remove it so the used
# range shrinks back down to A1:E3.
$ws.Range("F1:F3").Clear()

# Fix the header text: MODEL_CONDITION -> MODELCONDITION (now in column D
# after the shift above).
$ws.Range("D1").Value = "MODELCONDITION"
